$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Jon": bump every year in column A by one (2020..2054 -> 2021..2055)
# and drop what is now the extra trailing row (old row 37, year 2055).
# ---------------------------------------------------------------------------
$wsJon = $wb.Worksheets.Item("Jon")
for ($r = 2; $r -le 36; $r++) {
    $cell = $wsJon.Cells.Item($r, 1)
    $cur = $cell.Value2
    $cell.Value = $cur + 1
}
$wsJon.Rows.Item(37).Delete() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Jane": same treatment, one row longer (2020..2057 -> 2021..2058).
# ---------------------------------------------------------------------------
$wsJane = $wb.Worksheets.Item("Jane")
for ($r = 2; $r -le 39; $r++) {
    $cell = $wsJane.Cells.Item($r, 1)
    $cur = $cell.Value2
    $cell.Value = $cur + 1
}
$wsJane.Rows.Item(40).Delete() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Fixed Assets": insert a new "year" column before the old column D
# (basis), shifting basis/value/rate/yod/commission one column to the right.
# ---------------------------------------------------------------------------
$wsFixed = $wb.Worksheets.Item("Fixed Assets")
$wsFixed.Columns.Item(4).Insert() | Out-Null
$wsFixed.Range("D1").Value = "year"

# ---------------------------------------------------------------------------
# Restore per-sheet selections to match the saved workbook state.
# ---------------------------------------------------------------------------
$wsJon.Activate() | Out-Null
$wsJon.Range("A2:A36").Select() | Out-Null

$wsJane.Activate() | Out-Null
$wsJane.Range("B3").Select() | Out-Null

$wsDebts = $wb.Worksheets.Item("Debts")
$wsDebts.Activate() | Out-Null
$wsDebts.Range("D1:D1048576").Select() | Out-Null

# "Fixed Assets" is activated last so it ends up as the saved active tab.
$wsFixed.Activate() | Out-Null
$wsFixed.Range("D1:D1048576").Select() | Out-Null

Write-Output "done"
